# Applies the "Updated cryptos list" data refresh to the crypto tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds values such as "42.614.56" or "318.21" that are
# really text (thousand-separated prices), not numbers. Force the column to
# Text format first so Excel doesn't "helpfully" reinterpret single-decimal
# looking strings (e.g. "318.21") as floating point numbers and introduce
# binary rounding artifacts.
$ws.Columns("D").NumberFormat = "@"

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "42.614.56"
$ws.Range("E2").Value = "  -0.04%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "2.303.34"

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.14%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "318.21"
$ws.Range("E5").Value = "  -0.33%  "

# --- Row 6: Solana ---
$ws.Range("D6").Value = "103.78"
$ws.Range("E6").Value = "  +0.14%  "

# --- Row 7: XRP ---
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  -0.67%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  -0.21%  "

# --- Row 9: Cardano ---
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -1.50%  "

# --- Row 10: Avalanche ---
$ws.Range("D10").Value = "39.76"
$ws.Range("E10").Value = "  -0.14%  "

# --- Row 11: Dogecoin ---
$ws.Range("D11").Value = "0.0909"
$ws.Range("E11").Value = "  -0.40%  "

# --- Row 12: Polkadot ---
$ws.Range("D12").Value = "8.49"
$ws.Range("E12").Value = "  +1.54%  "

# --- Row 13: TRON ---
$ws.Range("E13").Value = "  +2.42%  "

# --- Row 14: Polygon ---
$ws.Range("E14").Value = "  +3.80%  "

# --- Row 15: Chainlink ---
$ws.Range("D15").Value = "15.38"
$ws.Range("E15").Value = "  +0.01%  "

# --- Row 16: WrappedliquidstakedEther2.0 ---
$ws.Range("D16").Value = "2.653.27"
$ws.Range("E16").Value = "  -0.27%  "

# --- Row 17: WrappedEther ---
$ws.Range("D17").Value = "2.301.90"
$ws.Range("E17").Value = "  -0.51%  "

# --- Row 18: WrappedBTC ---
$ws.Range("D18").Value = "42.701.67"
$ws.Range("E18").Value = "  +0.25%  "

# --- Rows 19 & 20: Uniswap / InternetComputer(DFINITY) swapped order ---
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").Value = "14.54"
$ws.Range("E19").Value = "  +33.77%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "7.53"
$ws.Range("E20").Value = "  +0.27%  "

# --- Row 21: ShibaInu ---
$ws.Range("E21").Value = "  +0.13%  "

# --- Row 22: Litecoin ---
$ws.Range("D22").Value = "73.94"
$ws.Range("E22").Value = "  +0.81%  "

# --- Row 23: PancakeSwap ---
$ws.Range("D23").Value = "3.55"
$ws.Range("E23").Value = "  -1.54%  "

# --- Row 24: BitcoinCash ---
$ws.Range("D24").Value = "266.63"
$ws.Range("E24").Value = "  -5.41%  "

# --- Row 25: ImmutableX ---
$ws.Range("E25").Value = "  -1.81%  "

# --- Row 26: Dai ---
$ws.Range("E26").Value = "  +0.59%  "

# --- Row 27: Cosmos ---
$ws.Range("D27").Value = "10.96"
$ws.Range("E27").Value = "  +0.54%  "

# --- Row 28: Toncoin ---
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  -3.49%  "

# --- Row 29: Filecoin ---
$ws.Range("D29").Value = "6.72"
$ws.Range("E29").Value = "  +13.36%  "

# --- Row 30: EthereumClassic ---
$ws.Range("D30").Value = "22.62"
$ws.Range("E30").Value = "  -1.46%  "

# --- Row 31: InjectiveProtocol ---
$ws.Range("D31").Value = "37.37"
$ws.Range("E31").Value = "  +3.25%  "

# --- Row 32: Monero ---
$ws.Range("D32").Value = "166.10"
$ws.Range("E32").Value = "  +0.82%  "

# --- Row 33: Hedera ---
$ws.Range("D33").Value = "0.0883"
$ws.Range("E33").Value = "  +0.71%  "

# --- Row 34: Stellar ---
$ws.Range("E34").Value = "  -3.48%  "

# --- Row 35: WEMIXToken ---
$ws.Range("D35").Value = "2.62"
$ws.Range("E35").Value = "  -0.19%  "

# --- Row 36: Kaspa ---
$ws.Range("E36").Value = "  -3.06%  "

# --- Row 37: RenderToken ---
$ws.Range("D37").Value = "4.56"
$ws.Range("E37").Value = "  -1.50%  "

# --- Row 38: VeChain ---
$ws.Range("D38").Value = "0.0354"
$ws.Range("E38").Value = "  -2.97%  "

# --- Row 39: NEARProtocol ---
$ws.Range("D39").Value = "3.72"
$ws.Range("E39").Value = "  -0.12%  "

# --- Row 40: LidoDAOToken ---
$ws.Range("D40").Value = "2.70"
$ws.Range("E40").Value = "  -3.19%  "

# --- Row 41: ARBITRUM ---
$ws.Range("E41").Value = "  +5.39%  "

# --- Row 42: MultiversX ---
$ws.Range("D42").Value = "70.44"
$ws.Range("E42").Value = "  +1.08%  "

# --- Row 43: BitcoinSV ---
$ws.Range("D43").Value = "95.79"
$ws.Range("E43").Value = "  -4.19%  "

# --- Row 44: Algorand ---
$ws.Range("E44").Value = "  +0.88%  "

# --- Row 45: FirstDigitalUSD ---
$ws.Range("E45").Value = "  +0.00%  "

# --- Row 46: Celestia ---
$ws.Range("D46").Value = "12.27"
$ws.Range("E46").Value = "  +0.43%  "

# --- Row 47: ordi ---
$ws.Range("D47").Value = "81.49"
$ws.Range("E47").Value = "  +2.43%  "

# --- Row 48: Aave ---
$ws.Range("D48").Value = "114.78"
$ws.Range("E48").Value = "  +1.53%  "

# --- Row 49: Maker ---
$ws.Range("D49").Value = "1.688.38"
$ws.Range("E49").Value = "  +4.45%  "

# --- Row 50: FraxShare ---
$ws.Range("D50").Value = "8.83"
$ws.Range("E50").Value = "  -1.51%  "

# --- Row 51: THORChain ---
$ws.Range("D51").Value = "5.18"
$ws.Range("E51").Value = "  -2.85%  "
